$wb = $excel.ActiveWorkbook

# Work on the "Repayment Schedule" sheet: insert a new (blank) column at N,
# shifting the old N:P ("Late", "Heading", "Outstanding") block to O:Q.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Update the selection on the Repayment Schedule sheet.
$wsSchedule.Range("S11").Select()

# Make "Repayment Schedule" the active sheet (was "Transactions").
$wsSchedule.Activate()
